$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company name (B column) reassignments ---
$ws.Range("B3").Value = 'IGI Holdings Limited (KASE:IGIHL)'
$ws.Range("B5").Value = 'IGI Life Insurance Limited (KASE:IGIL)'
$ws.Range("B6").Value = 'Askari Life Assurance Company Limited (KASE:ALAC)'
$ws.Range("B7").Value = 'Jubilee Life Insurance Company Limited (KASE:JLICL)'

# --- Clear cells removed in the update ---
$ws.Range("E5").Value = $null
$ws.Range("T6").Value = $null
$ws.Range("AO6").Value = $null
$ws.Range("AQ6").Value = $null

# --- Updated / new numeric values, grouped by row ---

# Row 2
$ws.Range("D2").Value = 0.09230000000000001
$ws.Range("E2").Value = 0.08460000000000001
$ws.Range("G2").Value = 0.05566983949885523
$ws.Range("H2").Value = 0.05566983949885523
$ws.Range("I2").Value = 0.05521187240084725
$ws.Range("J2").Value = 0.04538201279046201
$ws.Range("K2").Value = 32.267
$ws.Range("L2").Value = 0.03750101693339378
$ws.Range("M2").Value = 19.512
$ws.Range("N2").Value = 0.03267301863728462
$ws.Range("O2").Value = 0.6047044968543713
$ws.Range("P2").Value = 19.512
$ws.Range("Q2").Value = 0.03267301863728462
$ws.Range("R2").Value = 0.6047044968543713
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 173.718
$ws.Range("V2").Value = 0.2908923458195884
$ws.Range("W2").Value = 0.02323667711598746
$ws.Range("X2").Value = 0.1189149735414421
$ws.Range("Y2").Value = -0.09567829642545461
$ws.Range("Z2").Value = 6.432340325433815
$ws.Range("AA2").Value = -0.1031182973106748
$ws.Range("AB2").Value = 0.115953348461252
$ws.Range("AC2").Value = -0.218812658606795
$ws.Range("AD2").Value = 41.951
$ws.Range("AE2").Value = 0.1352431506950378
$ws.Range("AF2").Value = 42.08624315069504
$ws.Range("AG2").Value = -131.631756849305
$ws.Range("AH2").Value = 0.06583420485527749
$ws.Range("AI2").Value = 0.09373636170169057
$ws.Range("AJ2").Value = -0.2827396116079455
$ws.Range("AK2").Value = -0.4781962261206764
$ws.Range("AL2").Value = 4.108
$ws.Range("AM2").Value = 4.108
$ws.Range("AN2").Value = 0.794601761530448
$ws.Range("AO2").Value = 11.56158714703018
$ws.Range("AP2").Value = -2.493261802240837
$ws.Range("AQ2").Value = 11.56158714703018

# Row 3
$ws.Range("D3").Value = 0.09230000000000001
$ws.Range("E3").Value = -0.137
$ws.Range("G3").Value = 0.1283045977011494
$ws.Range("H3").Value = 0.1283045977011494
$ws.Range("I3").Value = 0.115948275862069
$ws.Range("J3").Value = 0.08148470890908337
$ws.Range("K3").Value = 5.93
$ws.Range("L3").Value = 0.08520114942528736
$ws.Range("M3").Value = 2.53
$ws.Range("N3").Value = 0.01396247240618102
$ws.Range("O3").Value = 0.4266441821247892
$ws.Range("P3").Value = 2.53
$ws.Range("Q3").Value = 0.01396247240618102
$ws.Range("R3").Value = 0.4266441821247892
$ws.Range("U3").Value = 5.61
$ws.Range("V3").Value = 0.03096026490066225
$ws.Range("W3").Value = 0.02323667711598746
$ws.Range("X3").Value = 0.1249474096164819
$ws.Range("Y3").Value = -0.1017107325004945
$ws.Range("Z3").Value = 0.251091309210289
$ws.Range("AA3").Value = 0.02046010224060104
$ws.Range("AB3").Value = 0.1165130538957476
$ws.Range("AC3").Value = -0.0960529516551465
$ws.Range("AD3").Value = 24.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 24.5
$ws.Range("AG3").Value = 18.89
$ws.Range("AH3").Value = 0.1191054934370442
$ws.Range("AI3").Value = 0.07860121912094964
$ws.Range("AJ3").Value = 0.09440751661752213
$ws.Range("AK3").Value = 0.06171387500408378
$ws.Range("AL3").Value = 2.84
$ws.Range("AM3").Value = 2.84
$ws.Range("AN3").Value = 2.740492170022371
$ws.Range("AO3").Value = 2.841549295774648
$ws.Range("AP3").Value = 2.112975391498882
$ws.Range("AQ3").Value = 2.841549295774648

# Row 4
$ws.Range("D4").Value = 0.0602
$ws.Range("E4").Value = 0.08460000000000001
$ws.Range("G4").Value = 0.05413105413105414
$ws.Range("H4").Value = 0.05413105413105414
$ws.Range("I4").Value = 0.05508072174738841
$ws.Range("J4").Value = 0.03865039281705948
$ws.Range("K4").Value = 12.3
$ws.Range("L4").Value = 0.03893637226970561
$ws.Range("M4").Value = 9.050000000000001
$ws.Range("N4").Value = 0.06929555895865239
$ws.Range("O4").Value = 0.7357723577235773
$ws.Range("P4").Value = 9.050000000000001
$ws.Range("Q4").Value = 0.06929555895865239
$ws.Range("R4").Value = 0.7357723577235773
$ws.Range("U4").Value = 135.9
$ws.Range("V4").Value = 1.040581929555896
$ws.Range("W4").Value = 0.3738601823708207
$ws.Range("X4").Value = 0.117250126394748
$ws.Range("Y4").Value = 0.2566100559760727
$ws.Range("Z4").Value = -2.244723939458537
$ws.Range("AA4").Value = -0.0867594620259297
$ws.Range("AB4").Value = 0.1158452679447433
$ws.Range("AC4").Value = -0.202604729970673
$ws.Range("AD4").Value = 3.11
$ws.Range("AF4").Value = 3.11
$ws.Range("AG4").Value = -132.79
$ws.Range("AH4").Value = 0.02325929249869119
$ws.Range("AI4").Value = 0.08335566872152238
$ws.Range("AJ4").Value = 60.63470319634709
$ws.Range("AK4").Value = 1.3468911654326
$ws.Range("AN4").Value = 0.1603092783505155
$ws.Range("AP4").Value = -6.844845360824742

# Row 5
$ws.Range("D5").Value = 0.08650000000000001
$ws.Range("G5").Value = -0.03028322440087146
$ws.Range("H5").Value = -0.03028322440087146
$ws.Range("I5").Value = -0.01361655773420479
$ws.Range("J5").Value = -0.01361655773420479
$ws.Range("K5").Value = -0.423
$ws.Range("L5").Value = -0.009215686274509804
$ws.Range("M5").Value = 0.002
$ws.Range("N5").Value = 0.00003273322422258593
$ws.Range("O5").Value = -0.004728132387706856
$ws.Range("P5").Value = 0.002
$ws.Range("Q5").Value = 0.00003273322422258593
$ws.Range("R5").Value = -0.004728132387706856
$ws.Range("U5").Value = 5.36
$ws.Range("V5").Value = 0.08772504091653029
$ws.Range("W5").Value = -0.03021428571428571
$ws.Range("X5").Value = 0.1164268478710277
$ws.Range("Y5").Value = -0.1466411335853134
$ws.Range("Z5").Value = 7.573007754495958
$ws.Range("AA5").Value = -0.1031182973106748
$ws.Range("AB5").Value = 0.1156943612961202
$ws.Range("AC5").Value = -0.218812658606795
$ws.Range("AD5").Value = 0.727
$ws.Range("AF5").Value = 0.727
$ws.Range("AG5").Value = -4.633
$ws.Range("AH5").Value = 0.01175861678554677
$ws.Range("AI5").Value = 0.05334996697732443
$ws.Range("AJ5").Value = -0.08204792179503073
$ws.Range("AK5").Value = -0.5604209507681143
$ws.Range("AL5").Value = 0.058
$ws.Range("AM5").Value = 0.058
$ws.Range("AN5").Value = -2.991769547325103
$ws.Range("AO5").Value = -10.77586206896552
$ws.Range("AP5").Value = 19.06584362139918
$ws.Range("AQ5").Value = -10.77586206896552

# Row 6
$ws.Range("D6").Value = 0.8270000000000001
$ws.Range("G6").Value = -0.8497409326424871
$ws.Range("H6").Value = -0.8497409326424871
$ws.Range("I6").Value = -0.7456210518854961
$ws.Range("J6").Value = -0.7456210518854961
$ws.Range("K6").Value = -1.44
$ws.Range("L6").Value = -0.7461139896373057
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 0.548
$ws.Range("V6").Value = 0.07517146776406036
$ws.Range("W6").Value = -0.6371681415929203
$ws.Range("X6").Value = 0.1189149735414421
$ws.Range("Y6").Value = -0.7560831151343624
$ws.Range("Z6").Value = 0.8595950952584128
$ws.Range("AA6").Value = -0.640932199122191
$ws.Range("AB6").Value = 0.115953348461252
$ws.Range("AC6").Value = -0.7568855475834431
$ws.Range("AD6").Value = 0.214
$ws.Range("AE6").Value = 0.1352431506950378
$ws.Range("AF6").Value = 0.3492431506950378
$ws.Range("AG6").Value = -0.1987568493049622
$ws.Range("AH6").Value = 0.04571698318874204
$ws.Range("AI6").Value = 0.3331699810902385
$ws.Range("AJ6").Value = -0.0280284916313272
$ws.Range("AK6").Value = -0.3973204811076563
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = -0.1643625192012289
$ws.Range("AP6").Value = 0.152655030188143

# Row 7
$ws.Range("D7").Value = 0.157
$ws.Range("E7").Value = 0.115
$ws.Range("G7").Value = 0.05830016389604308
$ws.Range("H7").Value = 0.05830016389604308
$ws.Range("I7").Value = 0.05642706626082884
$ws.Range("J7").Value = 0.03979989073597128
$ws.Range("K7").Value = 15.9
$ws.Range("L7").Value = 0.03722781549988293
$ws.Range("M7").Value = 7.93
$ws.Range("N7").Value = 0.03654377880184332
$ws.Range("O7").Value = 0.4987421383647798
$ws.Range("P7").Value = 7.93
$ws.Range("Q7").Value = 0.03654377880184332
$ws.Range("R7").Value = 0.4987421383647798
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 26.3
$ws.Range("V7").Value = 0.1211981566820277
$ws.Range("W7").Value = 0.235207100591716
$ws.Range("X7").Value = 0.1198715667467763
$ws.Range("Y7").Value = 0.1153355338449397
$ws.Range("Z7").Value = -38.82727272727268
$ws.Range("AA7").Value = -1.54532121212121
$ws.Range("AB7").Value = 0.1160482443602625
$ws.Range("AC7").Value = -1.661369456481473
$ws.Range("AD7").Value = 13.4
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 13.4
$ws.Range("AG7").Value = -12.9
$ws.Range("AH7").Value = 0.05815972222222222
$ws.Range("AI7").Value = 0.1570926143024619
$ws.Range("AJ7").Value = -0.06320431161195493
$ws.Range("AK7").Value = -0.2186440677966101
$ws.Range("AL7").Value = 1.21
$ws.Range("AM7").Value = 1.21
$ws.Range("AN7").Value = 0.5153846153846154
$ws.Range("AO7").Value = 19.91735537190083
$ws.Range("AP7").Value = -0.4961538461538462
$ws.Range("AQ7").Value = 19.91735537190083
